$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (A1:U1) shared-string values from *_old/*_new to *_FV2304/*_FV2310
$cols = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

$colIndex = 1
foreach ($c in $cols) {
    $ws.Cells.Item(1, $colIndex).Value = "$($c)_FV2304"
    $colIndex++
}
# column K (11) stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"
$colIndex = 12
foreach ($c in $cols) {
    $ws.Cells.Item(1, $colIndex).Value = "$($c)_FV2310"
    $colIndex++
}

# 2. Freeze the header row (split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into an Excel Table ("Table1")
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$tbl.Name = "Table1"
